# TMTI0046471 - Companies test data update (8th Oct 2025)
# Adds a new "System Admin" user (Ajay Nair) row to the Users sheet, and
# updates the active sheet / selection state so the Users sheet (not
# Contacts) is the one left active when the workbook is saved.

$wb = $excel.ActiveWorkbook

# --- Users sheet: append the new row ---------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A3").Value = "Ajay Nair"
$wsUsers.Range("B3").Value = "System Admin"

# --- Contacts sheet: it is no longer the active tab, record its last
#     selected cell ----------------------------------------------------
$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Range("G5").Select()

# --- Make Users the active sheet/tab and set its last selected cell --
$wsUsers.Activate()
$wsUsers.Range("B8").Select()
